$d = $word.ActiveDocument

# Update the title/date paragraph (unique text, single occurrence, safe to use ReplaceOne)
$d.Content.Find.Execute("2023-10-17 Tuesday", $true, $false, $false, $false, $false, $true, 0, $false, "2023-10-18 Wednesday", 1)

# Update table cells individually (row, col) -> new value, using the first table in the document.
$t = $d.Tables.Item(1)

$edits = @(
    @{Row=1;  Col=1; Old="36÷2=18, 0";  New="55÷3=18, 1"},
    @{Row=1;  Col=2; Old="63÷2=31, 1";  New="62÷7=8, 6"},
    @{Row=1;  Col=3; Old="89÷5=17, 4";  New="73÷6=12, 1"},
    @{Row=1;  Col=4; Old="20÷8=2, 4";   New="41÷6=6, 5"},
    @{Row=1;  Col=5; Old="52÷6=8, 4";   New="64÷7=9, 1"},

    @{Row=5;  Col=1; Old="25÷8=3, 1";   New="97÷6=16, 1"},
    @{Row=5;  Col=2; Old="44÷2=22, 0";  New="85÷3=28, 1"},
    @{Row=5;  Col=3; Old="86÷4=21, 2";  New="93÷9=10, 3"},
    @{Row=5;  Col=4; Old="90÷6=15, 0";  New="73÷7=10, 3"},
    @{Row=5;  Col=5; Old="95÷7=13, 4";  New="29÷2=14, 1"},

    @{Row=9;  Col=1; Old="82÷6=13, 4";  New="85÷9=9, 4"},
    @{Row=9;  Col=2; Old="33÷2=16, 1";  New="77÷9=8, 5"},
    @{Row=9;  Col=3; Old="45÷5=9, 0";   New="91÷2=45, 1"},
    @{Row=9;  Col=4; Old="86÷6=14, 2";  New="77÷9=8, 5"},
    @{Row=9;  Col=5; Old="86÷8=10, 6";  New="26÷5=5, 1"},

    @{Row=13; Col=1; Old="20÷2=10, 0";  New="43÷6=7, 1"},
    @{Row=13; Col=2; Old="26÷6=4, 2";   New="86÷8=10, 6"},
    @{Row=13; Col=3; Old="52÷5=10, 2";  New="82÷3=27, 1"},
    @{Row=13; Col=4; Old="52÷6=8, 4";   New="57÷9=6, 3"},
    @{Row=13; Col=5; Old="95÷2=47, 1";  New="28÷9=3, 1"},

    @{Row=17; Col=1; Old="15÷8=1, 7";   New="98÷8=12, 2"},
    @{Row=17; Col=2; Old="34÷5=6, 4";   New="21÷6=3, 3"},
    @{Row=17; Col=3; Old="97÷2=48, 1";  New="33÷4=8, 1"},
    @{Row=17; Col=4; Old="29÷8=3, 5";   New="84÷8=10, 4"},
    @{Row=17; Col=5; Old="36÷5=7, 1";   New="26÷8=3, 2"}
)

foreach ($e in $edits) {
    # Use wdFindStop (0) and wdReplaceOne (1) so the replace only touches the single
    # match inside this cell's Range, avoiding accidental document-wide replacement
    # of other cells that happen to contain the same (duplicate) text.
    $cellRange = $d.Tables.Item(1).Cell($e.Row, $e.Col).Range
    $cellRange.Find.Execute($e.Old, $true, $false, $false, $false, $false, $true, 0, $false, $e.New, 1)
}
